# Auto-generated edit script: updates cryptos price/volume data
# Forces each written cell to remain Text (matches source inlineStr
# cells) and resets the style afterwards so no stray number-format
# style gets attached to the cell (keeps s="0"/unset, like before).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '70.857.03'
Set-TextValue 'E2' '  +1.44%  '
Set-TextValue 'D3' '3.798.34'
Set-TextValue 'E3' '  +0.14%  '
Set-TextValue 'E4' '  +0.12%  '
Set-TextValue 'D5' '701.05'
Set-TextValue 'E5' '  +6.05%  '
Set-TextValue 'D6' '173.06'
Set-TextValue 'E6' '  +4.18%  '
Set-TextValue 'D7' '3.796.36'
Set-TextValue 'E7' '  +0.05%  '
Set-TextValue 'E8' '  +0.01%  '
Set-TextValue 'D9' '0.527'
Set-TextValue 'E9' '  -0.12%  '
Set-TextValue 'E10' '  +2.00%  '
Set-TextValue 'D11' '7.25'
Set-TextValue 'E11' '  +3.92%  '
Set-TextValue 'E12' '  +0.32%  '
Set-TextValue 'E13' '  +6.50%  '
Set-TextValue 'D14' '36.06'
Set-TextValue 'E14' '  +1.69%  '
Set-TextValue 'D15' '4.434.21'
Set-TextValue 'E15' '  +0.50%  '
Set-TextValue 'D16' '3.796.89'
Set-TextValue 'E16' '  +0.50%  '
Set-TextValue 'D17' '70.839.28'
Set-TextValue 'E17' '  +1.69%  '
Set-TextValue 'E18' '  -0.50%  '
Set-TextValue 'D19' '7.19'
Set-TextValue 'E19' '  +0.76%  '
Set-TextValue 'E20' '  +0.22%  '
Set-TextValue 'D21' '10.94'
Set-TextValue 'E21' '  +8.62%  '
Set-TextValue 'D22' '480.07'
Set-TextValue 'E22' '  +1.97%  '
Set-TextValue 'E23' '  -0.22%  '
Set-TextValue 'D24' '84.01'
Set-TextValue 'E24' '  +1.83%  '
Set-TextValue 'E25' '  -0.91%  '
Set-TextValue 'D26' '12.27'
Set-TextValue 'E26' '  +0.10%  '
Set-TextValue 'D27' '10.58'
Set-TextValue 'E27' '  +2.29%  '
Set-TextValue 'E28' '  +1.78%  '
Set-TextValue 'D29' '3.945.52'
Set-TextValue 'E29' '  +0.54%  '
Set-TextValue 'E30' '  -0.14%  '
Set-TextValue 'D31' '3.13'
Set-TextValue 'E31' '  +13.27%  '
Set-TextValue 'D32' '7.60'
Set-TextValue 'E32' '  +3.59%  '
Set-TextValue 'E33' '  -0.52%  '
Set-TextValue 'D34' '0.188'
Set-TextValue 'E34' '  +5.79%  '
Set-TextValue 'D35' '29.44'
Set-TextValue 'E35' '  +1.22%  '
Set-TextValue 'B36' 'Binance-PegBSC-USD'
Set-TextValue 'C36' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D36' '1.05'
Set-TextValue 'E36' '  +4.63%  '
Set-TextValue 'B37' 'Aptos'
Set-TextValue 'C37' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D37' '9.25'
Set-TextValue 'E37' '  +2.24%  '
Set-TextValue 'E38' '  +1.30%  '
Set-TextValue 'D39' '3.43'
Set-TextValue 'E39' '  +3.79%  '
Set-TextValue 'D40' '6.01'
Set-TextValue 'E40' '  +2.09%  '
Set-TextValue 'D41' '2.24'
Set-TextValue 'E41' '  +9.65%  '
Set-TextValue 'D42' '0.986'
Set-TextValue 'E42' '  +2.26%  '
Set-TextValue 'E43' '  +0.25%  '
Set-TextValue 'E44' '  +0.02%  '
Set-TextValue 'D45' '0.000318'
Set-TextValue 'E45' '  +17.33%  '
Set-TextValue 'D46' '164.75'
Set-TextValue 'E46' '  +3.83%  '
Set-TextValue 'D47' '48.83'
Set-TextValue 'E47' '  +1.78%  '
Set-TextValue 'D48' '44.74'
Set-TextValue 'E48' '  -1.61%  '
Set-TextValue 'B49' 'ONDO'
Set-TextValue 'C49' 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 'D49' '1.39'
Set-TextValue 'E49' '  -0.30%  '
Set-TextValue 'B50' 'TheGraph'
Set-TextValue 'C50' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D50' '0.301'
Set-TextValue 'E50' '  +0.28%  '
Set-TextValue 'D51' '413.92'
Set-TextValue 'E51' '  +7.20%  '
